# Replace the obsolete RED LED part (5975004407F / OEPS030006) with the
# new part (VFHR1116P-4C82A-TR / OEPS030038) on the BOM sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FED3_v7.3 BOM")

# Row 15 is the "RED" LED line (Reference: CHG0,L0,L2 / Value: RED).
# Only the OEPSPN, MPN and Description columns change; Reference, Value,
# Footprint and Type stay the same.
$ws.Range("D15").Value = "OEPS030038"
$ws.Range("E15").Value = "VFHR1116P-4C82A-TR"
$ws.Range("G15").Value = "DIODE: LED 0603 Red 635nm 20mA [VFHR1116P-4C82A-TR] [0603]"

# Leave the selection where the author left it after making the edit.
$ws.Range("A15").Select()
